$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'35.492.79"
$ws.Range("E2").Value = "  +0.50%  "
# Row 3
$ws.Range("D3").Value = "'1.897.25"
$ws.Range("E3").Value = "  -0.89%  "
# Row 4
$ws.Range("E4").Value = "  -0.66%  "
# Row 5
$ws.Range("D5").Value = "'247.65"
$ws.Range("E5").Value = "  -3.09%  "
# Row 6
$ws.Range("D6").Value = "'0.692"
$ws.Range("E6").Value = "  -4.24%  "
# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.71%  "
# Row 8
$ws.Range("D8").Value = "'43.94"
$ws.Range("E8").Value = "  +7.79%  "
# Row 9
$ws.Range("D9").Value = "'0.351"
$ws.Range("E9").Value = "  -5.19%  "
# Row 10
$ws.Range("D10").Value = "'0.0740"
$ws.Range("E10").Value = "  -3.38%  "
# Row 11
$ws.Range("D11").Value = "'0.0970"
$ws.Range("E11").Value = "  -1.89%  "
# Row 12
$ws.Range("D12").Value = "'13.07"
$ws.Range("E12").Value = "  +1.59%  "
# Row 13
$ws.Range("D13").Value = "'2.173.76"
$ws.Range("E13").Value = "  -0.78%  "
# Row 14
$ws.Range("D14").Value = "'0.721"
$ws.Range("E14").Value = "  -0.80%  "
# Row 15
$ws.Range("D15").Value = "'4.91"
$ws.Range("E15").Value = "  -1.33%  "
# Row 16
$ws.Range("D16").Value = "'1.882.60"
$ws.Range("E16").Value = "  -1.83%  "
# Row 17
$ws.Range("D17").Value = "'35.486.66"
$ws.Range("E17").Value = "  +0.50%  "
# Row 18
$ws.Range("D18").Value = "'73.56"
$ws.Range("E18").Value = "  -1.47%  "
# Row 19
$ws.Range("D19").Value = "'0.0₃0822"
$ws.Range("E19").Value = "  -3.93%  "
# Row 20
$ws.Range("D20").Value = "'247.71"
$ws.Range("E20").Value = "  +1.29%  "
# Row 21
$ws.Range("D21").Value = "'12.81"
$ws.Range("E21").Value = "  -2.34%  "
# Row 22
$ws.Range("E22").Value = "  -3.26%  "
# Row 23
$ws.Range("E23").Value = "  -0.76%  "
# Row 24
$ws.Range("D24").Value = "'2.53"
$ws.Range("E24").Value = "  +5.71%  "
# Row 25
$ws.Range("E25").Value = "  -10.28%  "
# Row 26
$ws.Range("D26").Value = "'165.72"
$ws.Range("E26").Value = "  -0.88%  "
# Row 27
$ws.Range("D27").Value = "'8.46"
$ws.Range("E27").Value = "  -2.46%  "
# Row 28
$ws.Range("D28").Value = "'18.35"
$ws.Range("E28").Value = "  -2.43%  "
# Row 29
$ws.Range("E29").Value = "  -4.33%  "
# Row 30
$ws.Range("D30").Value = "'4.128.41"
$ws.Range("E30").Value = "  -0.03%  "
# Row 31
$ws.Range("E31").Value = "  +7.97%  "
# Row 32
$ws.Range("D32").Value = "'4.24"
$ws.Range("E32").Value = "  -2.86%  "
# Row 33
$ws.Range("D33").Value = "'0.0579"
$ws.Range("E33").Value = "  -1.75%  "
# Row 34
$ws.Range("D34").Value = "'4.23"
$ws.Range("E34").Value = "  -0.45%  "
# Row 35
$ws.Range("E35").Value = "  -0.70%  "
# Row 36
$ws.Range("E36").Value = "  -6.19%  "
# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.62"
$ws.Range("E37").Value = "  -19.30%  "
# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.00"
$ws.Range("E38").Value = "  -2.28%  "
# Row 39
$ws.Range("D39").Value = "'17.18"
$ws.Range("E39").Value = "  +0.01%  "
# Row 40
$ws.Range("D40").Value = "'0.0676"
$ws.Range("E40").Value = "  +4.30%  "
# Row 41
$ws.Range("D41").Value = "'97.23"
$ws.Range("E41").Value = "  +0.29%  "
# Row 42
$ws.Range("E42").Value = "  -2.83%  "
# Row 43
$ws.Range("E43").Value = "  -3.34%  "
# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'1.287.05"
$ws.Range("E44").Value = "  -3.81%  "
# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.35"
$ws.Range("E45").Value = "  -3.44%  "
# Row 46
$ws.Range("E46").Value = "  +7.14%  "
# Row 47
$ws.Range("E47").Value = "  -1.09%  "
# Row 48
$ws.Range("D48").Value = "'2.75"
$ws.Range("E48").Value = "  -0.81%  "
# Row 49
$ws.Range("D49").Value = "'12.07"
$ws.Range("E49").Value = "  +1.51%  "
# Row 50
$ws.Range("D50").Value = "'6.38"
$ws.Range("E50").Value = "  -5.58%  "
# Row 51
$ws.Range("D51").Value = "'43.19"
$ws.Range("E51").Value = "  -4.54%  "
